$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.524.86"
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").Value = "3.122.76"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.68"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.77"
$ws.Range("E6").Value = "  -7.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.116.42"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.65"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.43"
$ws.Range("E13").Value = "  -6.13%  "
$ws.Range("E14").Value = "  -3.42%  "
$ws.Range("D15").Value = "3.623.63"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "63.486.57"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "3.123.02"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "510.17"
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.71"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.31"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.51"
$ws.Range("E24").Value = "  -3.63%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.51"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.36"
$ws.Range("E28").Value = "  -5.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -9.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.49"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("E32").Value = "  -7.69%  "
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "59.46"
$ws.Range("E34").Value = "  +11.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "532.16"
$ws.Range("E35").Value = "  -10.32%  "
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.28"
$ws.Range("E37").Value = "  -6.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0417"
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0801"
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.122"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.079.71"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -8.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.18"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.255"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("B45").Value = "CoreDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.72"
$ws.Range("E45").Value = "  +75.89%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.08"
$ws.Range("E47").Value = "  -7.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.10"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  -6.31%  "
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("D51").Value = "0.0₃0514"
$ws.Range("E51").Value = "  -7.04%  "
